$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for rows 2-12
# from serial date 45221 (2023-10-22) to serial date 45224 (2023-10-25)
$newDate = (Get-Date -Year 2023 -Month 10 -Day 25).Date

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
